# Region VII_ELECTRIFICATION.xlsx update
# - Inserts two new helper columns (AS, AT) ahead of the old "Status as of
#   July 4, 2025" column, which is pushed from AS to AU.
# - For rows 2-27 (the original on-grid ENERGIZATION rows) a new "ongrid"
#   marker is written into AS, and the old AR value is preserved in AT.
# - For a handful of specific rows the previous AR value is kept (old AR's
#   contents moved to AT) and a brand-new "BBM ..." batch label is written
#   into AR.
# - For every other data row, the old AR value is simply moved into AT and
#   AR is left blank.
# - The AS2:AS500 dropdown validation is moved to AU2:AU500.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1)
# ---------------------------------------------------------------------------

# Remember the current AS1 text ("Status as of July 4, 2025") before we
# overwrite the cell; it ends up in AU1.
$oldStatusHeader = $ws.Cells.Item(1, 45).Value2

# Give the two new header cells (AS1, AT1) the same style as AR1 (bold,
# bordered, centered) by copying AR1's formatting onto them.
$ws.Range("AR1").Copy()
$ws.Range("AS1").PasteSpecial(-4122)
$ws.Range("AT1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 45).Value2 = "Unnamed: 44"
$ws.Cells.Item(1, 46).Value2 = "Unnamed: 45"
$ws.Cells.Item(1, 47).Value2 = $oldStatusHeader

# ---------------------------------------------------------------------------
# 2. Data rows (rows 2-500)
# ---------------------------------------------------------------------------

# Rows that get a literal "ongrid" marker written into AS (their old AR
# value is preserved by moving it to AT).
$ongridRows = 2..27

# Rows where a brand-new "BBM ..." value replaces AR, while the old AR
# value is preserved by moving it to AT.
$bbmValues = @{
    245 = "BBM 2024 SOLAR"
    246 = "BBM 2024 SOLAR"
    248 = "BBM 2022"
    250 = "BBM 2024 SOLAR"
    254 = "BBM 2024 SOLAR"
    452 = "BBM 2023 UPGRADE"
    457 = "BBM 2024 UPGRADE"
    472 = "BBM 2025 UPGRADE"
    477 = "BBM 2024 UPGRADE"
    478 = "BBM 2023 UPGRADE"
    479 = "BBM 2023 UPGRADE"
    480 = "BBM 2023 UPGRADE"
    481 = "BBM 2024 ONGRID"
    482 = "BBM 2024 UPGRADE"
    485 = "BBM 2024 UPGRADE"
    486 = "BBM 2024 UPGRADE"
    487 = "BBM 2024 UPGRADE"
    488 = "BBM 2024 UPGRADE"
    489 = "BBM 2024 UPGRADE"
    490 = "BBM 2024 UPGRADE"
    491 = "BBM 2024 UPGRADE"
    492 = "BBM 2024 UPGRADE"
    494 = "BBM 2025 ONGRID"
    496 = "BBM 2024 ONGRID"
    497 = "BBM 2024 ONGRID"
    498 = "BBM 2024 UPGRADE"
    499 = "BBM 2024 UPGRADE"
    500 = "BBM 2024 UPGRADE"
}

for ($r = 2; $r -le 500; $r++) {
    $oldAR = $ws.Cells.Item($r, 44).Value2

    if ($bbmValues.ContainsKey($r)) {
        # Keep AR, but replace its contents with the new BBM label; stash
        # the previous AR value in AT.
        $ws.Cells.Item($r, 46).Value2 = $oldAR
        $ws.Cells.Item($r, 44).Value2 = $bbmValues[$r]
    }
    elseif ($ongridRows -contains $r) {
        # New "ongrid" marker goes into AS; old AR value moves to AT; AR
        # itself is cleared.
        $ws.Cells.Item($r, 46).Value2 = $oldAR
        $ws.Cells.Item($r, 45).Value2 = "ongrid"
        $ws.Cells.Item($r, 44).ClearContents()
    }
    else {
        # Plain move: AR -> AT, AR cleared.
        $ws.Cells.Item($r, 46).Value2 = $oldAR
        $ws.Cells.Item($r, 44).ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 3. Move the dropdown validation from AS2:AS500 to AU2:AU500
# ---------------------------------------------------------------------------

$ws.Range("AS2:AS500").Validation.Delete()
$ws.Range("AU2:AU500").Validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
$ws.Range("AU2:AU500").Validation.IgnoreBlank = $true
$ws.Range("AU2:AU500").Validation.InCellDropdown = $true
$ws.Range("AU2:AU500").Validation.ShowInput = $false
$ws.Range("AU2:AU500").Validation.ShowError = $false

Write-Host "Applied Region VII_ELECTRIFICATION column restructuring"
